# FaxAddressBookData.xlsx edit
#
# Semantic content change: the FaxLine value "10384" on the
# EditAddressBook and DeleteAddressBook sheets (cell B2 on both) is
# corrected to "1000". Both cells carry the "quote-prefixed text" style
# (the FaxLine column stores these codes as text, not numbers), so the
# replacement value is entered the same way a user would type it in
# Excel to force text storage: with a leading apostrophe.
#
# The rest of the workbook's view/selection state is also nudged to
# match where the file was left after this edit: a handful of sheets
# have a different last-selected cell, the AddressCreate sheet is
# scrolled over one column, and the tab that ends up active/selected
# moves from AddressBookQuery to CreateRecipient.

$wb = $excel.ActiveWorkbook

# --- the actual data correction -------------------------------------
$wsEditAB = $wb.Worksheets.Item("EditAddressBook")
$wsEditAB.Range("B2").Value = "'1000"

$wsDeleteAB = $wb.Worksheets.Item("DeleteAddressBook")
$wsDeleteAB.Range("B2").Value = "'1000"

# --- per-sheet selection / scroll state -------------------------------
$wsAddressCreate = $wb.Worksheets.Item("AddressCreate")
$wsAddressCreate.Activate()
$wsAddressCreate.Range("E7").Select()
$excel.ActiveWindow.ScrollColumn = 2

$wsEditAB.Activate()
$wsEditAB.Range("C2").Select()

$wsDeleteAB.Activate()
$wsDeleteAB.Range("A17").Select()

$wsAddressBookQuery = $wb.Worksheets.Item("AddressBookQuery")
$wsAddressBookQuery.Activate()
$wsAddressBookQuery.Range("A4").Select()

# CreateRecipient is selected last, so it ends up the active/selected tab.
$wsCreateRecipient = $wb.Worksheets.Item("CreateRecipient")
$wsCreateRecipient.Activate()
$wsCreateRecipient.Range("C2").Select()
